# Apply cell value updates per the target diff (Jogos_do_Dia_Betfair_Back_Lay_2026-01-07.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 2.4
$ws.Range("G2").Value = 2.6
$ws.Range("H2").Value = 3.55
$ws.Range("I2").Value = 3.9
$ws.Range("J2").Value = 2.84
$ws.Range("V2").Value = 1.35
$ws.Range("W2").Value = 1.62
# Row 3
$ws.Range("G3").Value = 3.65
$ws.Range("H3").Value = 2.08
$ws.Range("J3").Value = 3.75
$ws.Range("K3").Value = 4.2
$ws.Range("N3").Value = 4.7
$ws.Range("P3").Value = 2.3
$ws.Range("Q3").Value = 1.62
$ws.Range("S3").Value = 2.56
$ws.Range("T3").Value = 1.62
$ws.Range("U3").Value = 2.38
$ws.Range("V3").Value = 1.8
$ws.Range("W3").Value = 1.37
$ws.Range("AD3").Value = 12
$ws.Range("AI3").Value = 30
$ws.Range("AO3").Value = 12.5
# Row 4
$ws.Range("G4").Value = 5.6
$ws.Range("I4").Value = 1.96
$ws.Range("V4").Value = 2.04
# Row 5
$ws.Range("F5").Value = 3.2
$ws.Range("G5").Value = 3.3
$ws.Range("H5").Value = 1.81
$ws.Range("I5").Value = 2.02
$ws.Range("J5").Value = 4.3
$ws.Range("K5").Value = 6.8
$ws.Range("T5").Value = 1.05
$ws.Range("V5").Value = 1.98
$ws.Range("W5").Value = 1.43
$ws.Range("X5").Value = 100
$ws.Range("Y5").Value = 50
$ws.Range("Z5").Value = 36
$ws.Range("AA5").Value = 42
$ws.Range("AB5").Value = 70
$ws.Range("AC5").Value = 23
$ws.Range("AD5").Value = 18.5
$ws.Range("AE5").Value = 60
$ws.Range("AF5").Value = 65
$ws.Range("AG5").Value = 32
$ws.Range("AH5").Value = 55
$ws.Range("AI5").Value = 170
$ws.Range("AJ5").Value = 200
$ws.Range("AK5").Value = 180
$ws.Range("AL5").Value = 310
$ws.Range("AM5").Value = 740
$ws.Range("AN5").Value = 430
# Row 6
$ws.Range("G6").Value = 34
$ws.Range("H6").Value = 2.14
$ws.Range("I6").Value = 870
$ws.Range("J6").Value = 1.09
$ws.Range("K6").Value = 950
$ws.Range("V6").Value = 1.03
$ws.Range("W6").Value = 1.03
# Row 7
$ws.Range("J7").Value = 1.09
$ws.Range("V7").Value = 1.02
# Row 8
$ws.Range("F8").Value = 1.39
$ws.Range("G8").Value = 1.4
$ws.Range("J8").Value = 5.2
$ws.Range("K8").Value = 5.3
$ws.Range("N8").Value = 3.6
$ws.Range("O8").Value = 1.36
$ws.Range("P8").Value = 1.88
$ws.Range("R8").Value = 1.33
$ws.Range("T8").Value = 2.56
$ws.Range("U8").Value = 1.6
$ws.Range("W8").Value = 3.5
$ws.Range("AA8").Value = 620
$ws.Range("AE8").Value = 280
$ws.Range("AH8").Value = 40
$ws.Range("AK8").Value = 17
$ws.Range("AO8").Value = 480
# Row 9
$ws.Range("F9").Value = 3
$ws.Range("G9").Value = 3.05
$ws.Range("H9").Value = 2.64
$ws.Range("I9").Value = 2.66
$ws.Range("N9").Value = 3.7
$ws.Range("O9").Value = 1.35
$ws.Range("R9").Value = 1.35
$ws.Range("S9").Value = 3.7
$ws.Range("U9").Value = 2.18
$ws.Range("V9").Value = 1.6
$ws.Range("W9").Value = 1.48
$ws.Range("Y9").Value = 10.5
$ws.Range("Z9").Value = 16.5
$ws.Range("AD9").Value = 11.5
$ws.Range("AE9").Value = 28
$ws.Range("AH9").Value = 17
$ws.Range("AK9").Value = 36
$ws.Range("AL9").Value = 48
# Row 10
$ws.Range("F10").Value = 3.55
$ws.Range("Q10").Value = 1.86
$ws.Range("S10").Value = 3.2
# Row 11
$ws.Range("F11").Value = 1.89
$ws.Range("G11").Value = 1.9
$ws.Range("H11").Value = 4.8
$ws.Range("I11").Value = 4.9
$ws.Range("P11").Value = 1.88
$ws.Range("R11").Value = 1.32
$ws.Range("W11").Value = 2.1
$ws.Range("Y11").Value = 15.5
$ws.Range("AE11").Value = 65
# Row 12
$ws.Range("F12").Value = 3.3
$ws.Range("H12").Value = 2.4
$ws.Range("I12").Value = 2.42
$ws.Range("M12").Value = 1.07
$ws.Range("Q12").Value = 1.98
$ws.Range("R12").Value = 1.38
$ws.Range("T12").Value = 1.76
$ws.Range("V12").Value = 1.7
$ws.Range("W12").Value = 1.42
$ws.Range("X12").Value = 14
$ws.Range("Y12").Value = 10.5
# Row 13
$ws.Range("F13").Value = 2.24
$ws.Range("J13").Value = 3.7
$ws.Range("K13").Value = 3.75
$ws.Range("N13").Value = 4.4
$ws.Range("O13").Value = 1.28
$ws.Range("Q13").Value = 1.86
$ws.Range("R13").Value = 1.45
$ws.Range("S13").Value = 3.15
# Row 14
$ws.Range("F14").Value = 1.43
$ws.Range("G14").Value = 1.44
$ws.Range("H14").Value = 8.199999999999999
$ws.Range("I14").Value = 8.4
$ws.Range("P14").Value = 3.2
$ws.Range("Q14").Value = 1.44
$ws.Range("T14").Value = 1.66
$ws.Range("V14").Value = 1.13
$ws.Range("W14").Value = 3.25
$ws.Range("X14").Value = 34
$ws.Range("Z14").Value = 80
$ws.Range("AJ14").Value = 14
$ws.Range("AL14").Value = 23
$ws.Range("AO14").Value = 70
# Row 15
$ws.Range("F15").Value = 1.86
$ws.Range("G15").Value = 1.87
$ws.Range("H15").Value = 5.2
$ws.Range("O15").Value = 1.38
$ws.Range("R15").Value = 1.3
$ws.Range("T15").Value = 2
$ws.Range("U15").Value = 1.96
$ws.Range("AD15").Value = 19.5
$ws.Range("AF15").Value = 10.5
$ws.Range("AK15").Value = 20
$ws.Range("AL15").Value = 40
$ws.Range("AN15").Value = 14
$ws.Range("AO15").Value = 100
# Row 16
$ws.Range("G16").Value = 2.6
$ws.Range("H16").Value = 3.35
$ws.Range("J16").Value = 3.1
$ws.Range("K16").Value = 3.15
$ws.Range("N16").Value = 2.96
$ws.Range("O16").Value = 1.5
$ws.Range("S16").Value = 4.8
$ws.Range("U16").Value = 1.91
$ws.Range("W16").Value = 1.62
$ws.Range("Y16").Value = 10.5
$ws.Range("AJ16").Value = 36
# Row 17
$ws.Range("G17").Value = 2.42
$ws.Range("J17").Value = 3.15
$ws.Range("P17").Value = 1.67
$ws.Range("Q17").Value = 2.42
$ws.Range("R17").Value = 1.25
$ws.Range("W17").Value = 1.7
$ws.Range("AA17").Value = 75
$ws.Range("AM17").Value = 140
# Row 18
$ws.Range("F18").Value = 9.800000000000001
$ws.Range("L18").Value = 1.32
$ws.Range("N18").Value = 4.9
$ws.Range("O18").Value = 1.24
$ws.Range("P18").Value = 2.34
$ws.Range("Q18").Value = 1.72
$ws.Range("T18").Value = 2
$ws.Range("U18").Value = 1.95
$ws.Range("V18").Value = 3.5
$ws.Range("Z18").Value = 8.4
$ws.Range("AF18").Value = 90
# Row 19
$ws.Range("H19").Value = 1.73
$ws.Range("J19").Value = 4.1
$ws.Range("K19").Value = 4.2
$ws.Range("P19").Value = 2.2
$ws.Range("Q19").Value = 1.8
$ws.Range("R19").Value = 1.46
$ws.Range("AA19").Value = 17
$ws.Range("AD19").Value = 9.800000000000001
$ws.Range("AG19").Value = 20
$ws.Range("AN19").Value = 75
$ws.Range("AO19").Value = 9
# Row 20
$ws.Range("F20").Value = 1.76
$ws.Range("I20").Value = 5.5
$ws.Range("T20").Value = 1.85
$ws.Range("U20").Value = 2.14
$ws.Range("Z20").Value = 40
$ws.Range("AC20").Value = 8.800000000000001
$ws.Range("AE20").Value = 65
$ws.Range("AJ20").Value = 17
$ws.Range("AM20").Value = 85
$ws.Range("AN20").Value = 9.800000000000001
